$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""64.095.95"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Formula = "=""3.513.79"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Formula = "=""586.44"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -2.83%  "

$ws.Range("D6").Formula = "=""132.57"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -3.30%  "

$ws.Range("D7").Formula = "=""3.512.32"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -1.88%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").Formula = "=""7.14"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").Formula = "=""0.387"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Formula = "=""4.111.93"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -1.80%  "

$ws.Range("D14").Formula = "=""27.88"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").Formula = "=""0.0000181"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Formula = "=""3.516.29"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Formula = "=""64.143.46"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -2.44%  "

$ws.Range("D19").Formula = "=""10.01"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Formula = "=""14.49"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("E21").Value = "  -3.36%  "

$ws.Range("D22").Formula = "=""393.60"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").Formula = "=""0.580"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Formula = "=""3.658.84"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -1.73%  "

$ws.Range("D25").Formula = "=""73.07"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Formula = "=""0.0000113"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -3.54%  "

$ws.Range("E28").Value = "  -1.57%  "

$ws.Range("E29").Value = "  -7.26%  "

$ws.Range("D30").Formula = "=""0.999"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Formula = "=""2.27"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("D32").Formula = "=""8.29"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -3.17%  "

$ws.Range("D33").Formula = "=""3.517.08"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Formula = "=""24.02"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").Formula = "=""0.145"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("D37").Formula = "=""5.39"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("D38").Formula = "=""7.02"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  -2.65%  "

$ws.Range("D40").Formula = "=""168.36"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").Formula = "=""0.0813"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.58%  "

$ws.Range("D42").Formula = "=""26.76"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("D43").Formula = "=""0.814"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.77%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Formula = "=""41.89"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("E46").Value = "  -5.25%  "

$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("D48").Formula = "=""1.66"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -2.69%  "

$ws.Range("D49").Formula = "=""2.455.25"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("D50").Formula = "=""6.90"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").Formula = "=""0.902"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.28%  "

$excel.CutCopyMode = 0